$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to be stored as text even when its content
# looks like a number (e.g. "1.004", "0.06140") so Excel doesn't
# coerce it into a numeric value and strip significant trailing
# zeros / reformat the decimal separators used by this price feed.
function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range("D2").Value = '27.698.27'
$ws.Range("E2").Value = '  -1.97%  '
$ws.Range("D3").Value = '1.756.90'
$ws.Range("E3").Value = '  -2.21%  '
Set-TextValue "D4" '1.004'
$ws.Range("E4").Value = '  +0.02%  '
Set-TextValue "D5" '325.36'
$ws.Range("E5").Value = '  -3.88%  '
Set-TextValue "D6" '1.001'
$ws.Range("E6").Value = '  +0.05%  '
Set-TextValue "D7" '0.4525'
$ws.Range("E7").Value = '  -1.35%  '
Set-TextValue "D8" '0.3728'
$ws.Range("E8").Value = '  -0.67%  '
Set-TextValue "D9" '45.38'
$ws.Range("E9").Value = '  +0.48%  '
Set-TextValue "D10" '0.07544'
$ws.Range("E10").Value = '  -0.60%  '
Set-TextValue "D11" '1.127'
$ws.Range("E11").Value = '  -1.54%  '
Set-TextValue "D12" '1.001'
$ws.Range("E12").Value = '  -0.10%  '
Set-TextValue "D13" '21.77'
$ws.Range("E13").Value = '  -2.35%  '
Set-TextValue "D14" '6.208'
$ws.Range("E14").Value = '  -1.34%  '
Set-TextValue "D15" '7.349'
$ws.Range("E15").Value = '  -1.91%  '
$ws.Range("D16").Value = '1.755.85'
$ws.Range("E16").Value = '  -2.37%  '
$ws.Range("E17").Value = '  -1.37%  '
Set-TextValue "D18" '87.96'
$ws.Range("E18").Value = '  +8.33%  '
Set-TextValue "D19" '0.06225'
$ws.Range("E19").Value = '  -7.47%  '
Set-TextValue "D20" '1.001'
$ws.Range("E20").Value = '  +0.02%  '
Set-TextValue "D21" '17.25'
$ws.Range("E21").Value = '  -1.12%  '
Set-TextValue "D22" '6.190'
$ws.Range("E22").Value = '  -2.82%  '
Set-TextValue "D23" '0.5318'
$ws.Range("E23").Value = '  -2.84%  '
$ws.Range("D24").Value = '27.749.27'
$ws.Range("E24").Value = '  -1.82%  '
Set-TextValue "D25" '11.67'
$ws.Range("E25").Value = '  -1.52%  '
Set-TextValue "D26" '2.319'
$ws.Range("E26").Value = '  -4.25%  '
Set-TextValue "D27" '20.66'
$ws.Range("E27").Value = '  -0.12%  '
Set-TextValue "D28" '153.29'
$ws.Range("E28").Value = '  +0.90%  '
Set-TextValue "D29" '2.357'
$ws.Range("E29").Value = '  +0.44%  '
$ws.Range("D30").Value = '1.954.63'
$ws.Range("E30").Value = '  -2.41%  '
Set-TextValue "D31" '128.44'
$ws.Range("E31").Value = '  -3.01%  '
Set-TextValue "D32" '1.219'
$ws.Range("E32").Value = '  -1.19%  '
Set-TextValue "D33" '0.09306'
$ws.Range("E33").Value = '  -2.29%  '
Set-TextValue "D34" '5.730'
$ws.Range("E34").Value = '  -1.20%  '
Set-TextValue "D35" '3.655'
$ws.Range("E35").Value = '  -9.18%  '
Set-TextValue "D36" '12.64'
$ws.Range("E36").Value = '  +5.25%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D37" '0.02329'
$ws.Range("E37").Value = '  -0.66%  '
$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue "D38" '0.2173'
$ws.Range("E38").Value = '  -6.33%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D39" '0.06140'
$ws.Range("E39").Value = '  -2.76%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue "D40" '0.6491'
$ws.Range("E40").Value = '  -1.52%  '
Set-TextValue "D41" '5.090'
$ws.Range("E41").Value = '  -3.01%  '
Set-TextValue "D42" '1.198'
$ws.Range("E42").Value = '  -2.69%  '
Set-TextValue "D43" '7.984'
$ws.Range("E43").Value = '  -4.74%  '
Set-TextValue "D44" '1.419'
$ws.Range("E44").Value = '  -4.29%  '
Set-TextValue "D45" '1.001'
$ws.Range("E45").Value = '  +0.02%  '
Set-TextValue "D46" '13.89'
$ws.Range("E46").Value = '  -2.63%  '
Set-TextValue "D47" '0.5988'
$ws.Range("E47").Value = '  -1.75%  '
Set-TextValue "D48" '3.755'
$ws.Range("E48").Value = '  -2.52%  '
Set-TextValue "D49" '126.17'
$ws.Range("E49").Value = '  -3.20%  '
Set-TextValue "D50" '1.991'
$ws.Range("E50").Value = '  -1.87%  '
Set-TextValue "D51" '0.06918'
$ws.Range("E51").Value = '  -3.20%  '
